$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a literal text value into a cell, even when the text looks
# like a number (e.g. "0.999"), without leaving any NumberFormat/style change
# behind. We flip the cell to text ("@") just long enough for the assignment
# to be stored as a string, then clear the formatting delta straight back off.
function Set-TextCell([object]$Range, [string]$Text, [bool]$NumericLooking) {
    if ($NumericLooking) {
        $Range.NumberFormat = "@"
        $Range.Value = $Text
        $Range.ClearFormats()
    } else {
        $Range.Value = $Text
    }
}

# Row 2
Set-TextCell $ws.Range("D2") "61.773.02" $false
Set-TextCell $ws.Range("E2") "  -1.62%  " $false

# Row 3
Set-TextCell $ws.Range("D3") "2.899.77" $false
Set-TextCell $ws.Range("E3") "  -2.38%  " $false

# Row 4
Set-TextCell $ws.Range("D4") "0.999" $true
Set-TextCell $ws.Range("E4") "  -0.03%  " $false

# Row 5
Set-TextCell $ws.Range("D5") "575.83" $true
Set-TextCell $ws.Range("E5") "  -3.23%  " $false

# Row 6
Set-TextCell $ws.Range("D6") "144.65" $true
Set-TextCell $ws.Range("E6") "  -0.78%  " $false

# Row 7
Set-TextCell $ws.Range("E7") "  +0.15%  " $false

# Row 8
Set-TextCell $ws.Range("D8") "0.504" $true
Set-TextCell $ws.Range("E8") "  -0.21%  " $false

# Row 9
Set-TextCell $ws.Range("D9") "2.899.36" $false
Set-TextCell $ws.Range("E9") "  -2.46%  " $false

# Row 10
Set-TextCell $ws.Range("D10") "6.67" $true
Set-TextCell $ws.Range("E10") "  -8.30%  " $false

# Row 11
Set-TextCell $ws.Range("D11") "0.150" $true
Set-TextCell $ws.Range("E11") "  +1.89%  " $false

# Row 12
Set-TextCell $ws.Range("D12") "0.432" $true
Set-TextCell $ws.Range("E12") "  -3.06%  " $false

# Row 13
Set-TextCell $ws.Range("D13") "0.0000235" $true
Set-TextCell $ws.Range("E13") "  -1.81%  " $false

# Row 14
Set-TextCell $ws.Range("D14") "32.55" $true
Set-TextCell $ws.Range("E14") "  -1.82%  " $false

# Row 15
Set-TextCell $ws.Range("E15") "  -0.75%  " $false

# Row 16
Set-TextCell $ws.Range("D16") "3.378.24" $false
Set-TextCell $ws.Range("E16") "  -2.41%  " $false

# Row 17
Set-TextCell $ws.Range("D17") "61.737.31" $false
Set-TextCell $ws.Range("E17") "  -1.51%  " $false

# Row 18
Set-TextCell $ws.Range("D18") "6.63" $true
Set-TextCell $ws.Range("E18") "  -1.17%  " $false

# Row 19
Set-TextCell $ws.Range("D19") "2.936.04" $false
Set-TextCell $ws.Range("E19") "  -0.68%  " $false

# Row 20
Set-TextCell $ws.Range("D20") "434.80" $true
Set-TextCell $ws.Range("E20") "  -1.59%  " $false

# Row 21
Set-TextCell $ws.Range("D21") "13.30" $true
Set-TextCell $ws.Range("E21") "  -1.41%  " $false

# Row 22
Set-TextCell $ws.Range("D22") "0.657" $true
Set-TextCell $ws.Range("E22") "  -2.18%  " $false

# Row 23
Set-TextCell $ws.Range("D23") "6.93" $true
Set-TextCell $ws.Range("E23") "  -1.98%  " $false

# Row 24
Set-TextCell $ws.Range("D24") "79.72" $true
Set-TextCell $ws.Range("E24") "  -2.37%  " $false

# Row 25
Set-TextCell $ws.Range("D25") "11.90" $true
Set-TextCell $ws.Range("E25") "  +0.08%  " $false

# Row 26
Set-TextCell $ws.Range("D26") "10.15" $true
Set-TextCell $ws.Range("E26") "  -9.91%  " $false

# Row 27
Set-TextCell $ws.Range("E27") "  -0.03%  " $false

# Row 28
Set-TextCell $ws.Range("D28") "2.04" $true
Set-TextCell $ws.Range("E28") "  -4.23%  " $false

# Row 29
Set-TextCell $ws.Range("D29") "0.0000108" $true
Set-TextCell $ws.Range("E29") "  +13.33%  " $false

# Row 30
Set-TextCell $ws.Range("D30") "7.03" $true
Set-TextCell $ws.Range("E30") "  -3.14%  " $false

# Row 31
Set-TextCell $ws.Range("D31") "2.53" $true
Set-TextCell $ws.Range("E31") "  -3.02%  " $false

# Row 32
Set-TextCell $ws.Range("D32") "2.09" $true
Set-TextCell $ws.Range("E32") "  -2.64%  " $false

# Row 35
Set-TextCell $ws.Range("D35") "25.70" $true
Set-TextCell $ws.Range("E35") "  -3.07%  " $false

# Row 36
Set-TextCell $ws.Range("D36") "0.962" $true
Set-TextCell $ws.Range("E36") "  -3.20%  " $false

# Row 37
Set-TextCell $ws.Range("D37") "3.03" $true
Set-TextCell $ws.Range("E37") "  -0.21%  " $false

# Row 38
Set-TextCell $ws.Range("D38") "5.45" $true
Set-TextCell $ws.Range("E38") "  -3.21%  " $false

# Row 39
Set-TextCell $ws.Range("D39") "49.09" $true
Set-TextCell $ws.Range("E39") "  -0.81%  " $false

# Row 40
Set-TextCell $ws.Range("D40") "1.96" $true
Set-TextCell $ws.Range("E40") "  -4.51%  " $false

# Row 41
Set-TextCell $ws.Range("E41") "  -1.71%  " $false

# Row 42
Set-TextCell $ws.Range("D42") "8.28" $true
Set-TextCell $ws.Range("E42") "  -3.13%  " $false

# Row 43
Set-TextCell $ws.Range("D43") "0.269" $true
Set-TextCell $ws.Range("E43") "  -4.67%  " $false

# Row 44
Set-TextCell $ws.Range("D44") "38.47" $true
Set-TextCell $ws.Range("E44") "  -5.02%  " $false

# Row 45
Set-TextCell $ws.Range("D45") "2.682.88" $false
Set-TextCell $ws.Range("E45") "  -2.27%  " $false

# Row 46
Set-TextCell $ws.Range("D46") "134.26" $true

# Row 47
Set-TextCell $ws.Range("D47") "0.0335" $true
Set-TextCell $ws.Range("E47") "  -1.65%  " $false

# Row 48
Set-TextCell $ws.Range("D48") "341.94" $true
Set-TextCell $ws.Range("E48") "  -6.05%  " $false

# Row 49
Set-TextCell $ws.Range("E49") "  +0.01%  " $false

# Row 50
Set-TextCell $ws.Range("E50") "  -1.85%  " $false

# Row 51
Set-TextCell $ws.Range("D51") "21.79" $true
Set-TextCell $ws.Range("E51") "  -5.54%  " $false

# Row 33 (was Hedera) becomes FirstDigitalUSD; row 34 (was FirstDigitalUSD) becomes Hedera.
$ws.Range("B33").Value = "FirstDigitalUSD"
$ws.Range("C33").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextCell $ws.Range("D33") "1.00" $true
Set-TextCell $ws.Range("E33") "  +0.37%  " $false

$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextCell $ws.Range("D34") "0.107" $true
Set-TextCell $ws.Range("E34") "  -3.19%  " $false
